# Rolling back Sources-And-Credits to merge with develop
# - Remove the extra freesound.org SFX credit rows (old rows 15-22)
# - Re-point hyperlinks so they line up with the remaining rows again
# - Restore selection to A15 (post-deletion)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing hyperlinks first - row deletion below does not
# automatically keep hyperlink anchors in sync, so we rebuild them
# from scratch afterwards.
$ws.Hyperlinks.Delete()

# Remove the 8 rows of now-unwanted freesound.org sound-effect credits
# (old rows 15-22: qubodup whoosh, schots gun-shot, rkkaleikau laser,
# Merrick079 punch2, RICHERlandTV heavy-impacts, sharesynth jump01,
# MattRuthSound punch, dersuperanton game-over). This shifts the old
# rows 24-26 up to become rows 16-18.
$ws.Rows("15:22").Delete()

# Re-create the hyperlinks that remain, in the same order/mapping as
# the target workbook so relationship ids line up (rId1..rId9).
$ws.Hyperlinks.Add($ws.Range("B6"), "https://opengameart.org/content/spikes-0")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://www.fontspace.com/a-area-kilometer-50-font-f53888")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://freesound.org/people/Whiprealgood/sounds/87535/")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://freesound.org/people/suntemple/sounds/253172/")
$ws.Hyperlinks.Add($ws.Range("B5"), "https://opengameart.org/content/simple-explosion-bleeds-game-art")
$ws.Hyperlinks.Add($ws.Range("B7"), "https://opengameart.org/content/various-inventory-24-pixel-icon-set")
$ws.Hyperlinks.Add($ws.Range("B17"), "https://elthen.itch.io/2d-pixel-art-vegetable-monsters-sprite-pack")
$ws.Hyperlinks.Add($ws.Range("B18"), "https://free-game-assets.itch.io/night-city-street-2d-background-tiles")
$ws.Hyperlinks.Add($ws.Range("B8"), "https://opengameart.org/content/energy-icon")

# Adding a hyperlink can register a fresh duplicate "Hyperlink" cell
# style instead of reusing the existing one; re-applying the named
# style per-cell collapses it back to the original style record.
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B7").Style = "Hyperlink"
$ws.Range("B8").Style = "Hyperlink"
$ws.Range("B17").Style = "Hyperlink"
$ws.Range("B18").Style = "Hyperlink"

# Restore the selection/active cell to A15, matching the saved view.
$ws.Range("A15").Select()
